# Update the K column (column G) values for rows 2-10 and 12 of Sheet1.
# The data was regenerated to use K (strikeouts) instead of Strike# counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 1
    10 = 2
    12 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
